$d = $word.ActiveDocument

$replacements = @(
    @("2025-05-23 Friday", "2025-05-24 Saturday"),
    @("432×7=3024", "669×9=6021"),
    @("204×3=612", "509×2=1018"),
    @("780×3=2340", "929×2=1858"),
    @("606×3=1818", "467×8=3736"),
    @("702×3=2106", "265×8=2120"),
    @("116×7=812", "254×9=2286"),
    @("495×9=4455", "176×3=528"),
    @("734×7=5138", "992×2=1984"),
    @("977×5=4885", "605×9=5445"),
    @("678×5=3390", "841×8=6728"),
    @("755×5=3775", "564×6=3384"),
    @("851×5=4255", "985×3=2955"),
    @("259×7=1813", "322×6=1932"),
    @("911×7=6377", "709×4=2836"),
    @("547×3=1641", "491×5=2455"),
    @("313×4=1252", "956×8=7648"),
    @("611×4=2444", "474×5=2370"),
    @("221×8=1768", "890×9=8010"),
    @("110×3=330", "225×2=450"),
    @("384×9=3456", "724×6=4344"),
    @("471×8=3768", "672×8=5376"),
    @("612×8=4896", "968×4=3872"),
    @("463×8=3704", "507×3=1521"),
    @("362×9=3258", "843×9=7587"),
    @("893×4=3572", "981×6=5886")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
